$wb = $excel.ActiveWorkbook

# ---- Sheet1: main results table (instances 1-10); MP time limit & recourse data fix ----
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("B2").Value = -963.1075412558371
$ws1.Range("C2").Value = 438.513680963
$ws1.Range("F2").Value = 100
$ws1.Range("G2").Value = 110400
$ws1.Range("H2").Value = 120500
$ws1.Range("I2").Value = 10000

$ws1.Range("B3").Value = -960.465247165562
$ws1.Range("C3").Value = 406.883121483
$ws1.Range("F3").Value = 100
$ws1.Range("G3").Value = 110400
$ws1.Range("H3").Value = 120500
$ws1.Range("I3").Value = 10000

$ws1.Range("B4").Value = -962.1299820331453
$ws1.Range("C4").Value = 1148.367505049
$ws1.Range("F4").Value = 100
$ws1.Range("G4").Value = 110400
$ws1.Range("H4").Value = 120500
$ws1.Range("I4").Value = 10000

$ws1.Range("B5").Value = -958.1263449022533
$ws1.Range("C5").Value = 831.356580329
$ws1.Range("F5").Value = 100
$ws1.Range("G5").Value = 110400
$ws1.Range("H5").Value = 120500
$ws1.Range("I5").Value = 10000

$ws1.Range("B6").Value = -963.1399746398824
$ws1.Range("C6").Value = 1021.811989351
$ws1.Range("F6").Value = 100
$ws1.Range("G6").Value = 110400
$ws1.Range("H6").Value = 120500
$ws1.Range("I6").Value = 10000

$ws1.Range("B7").Value = -953.069846438852
$ws1.Range("C7").Value = 516.375281886
$ws1.Range("F7").Value = 100
$ws1.Range("G7").Value = 110400
$ws1.Range("H7").Value = 120500
$ws1.Range("I7").Value = 10000

$ws1.Range("B8").Value = -965.4153674106826
$ws1.Range("C8").Value = 954.205450714
$ws1.Range("F8").Value = 100
$ws1.Range("G8").Value = 110400
$ws1.Range("H8").Value = 120500
$ws1.Range("I8").Value = 10000

$ws1.Range("B9").Value = -954.6051029237938
$ws1.Range("C9").Value = 1215.461226395
$ws1.Range("F9").Value = 100
$ws1.Range("G9").Value = 110400
$ws1.Range("H9").Value = 120500
$ws1.Range("I9").Value = 10000

$ws1.Range("B10").Value = -960.9871769768399
$ws1.Range("C10").Value = 422.440112962
$ws1.Range("F10").Value = 100
$ws1.Range("G10").Value = 110400
$ws1.Range("H10").Value = 120500
$ws1.Range("I10").Value = 10000

$ws1.Range("B11").Value = -965.36122253593
$ws1.Range("C11").Value = 953.306023174
$ws1.Range("F11").Value = 100
$ws1.Range("G11").Value = 110400
$ws1.Range("H11").Value = 120500
$ws1.Range("I11").Value = 10000

# ---- Per-instance detail sheets (named "1".."10"): MP time limit increase + recourse fix ----

$ws = $wb.Worksheets.Item("1")
$ws.Range("D2").Value = 0.8462853815408935
$ws.Range("E2").Value = 136.22463
$ws.Range("B3").Value = -963.1075412558371
$ws.Range("C3").Value = 0.002352205422191391
$ws.Range("D3").Value = 425.3277832817665

$ws = $wb.Worksheets.Item("2")
$ws.Range("D2").Value = 0.07757260252770996
$ws.Range("E2").Value = 137.1635
$ws.Range("B3").Value = -960.465247165562
$ws.Range("C3").Value = 0.0008404312086639922
$ws.Range("D3").Value = 400.785045977442

$ws = $wb.Worksheets.Item("3")
$ws.Range("D2").Value = 0.13341643966540528
$ws.Range("E2").Value = 136.59474
$ws.Range("B3").Value = -962.1299820331453
$ws.Range("C3").Value = 0.05491416457396366
$ws.Range("D3").Value = 1142.279286326067

$ws = $wb.Worksheets.Item("4")
$ws.Range("D2").Value = 0.10627470621508789
$ws.Range("E2").Value = 137.84913
$ws.Range("B3").Value = -958.1263449022533
$ws.Range("C3").Value = 0.09879981259225552
$ws.Range("D3").Value = 824.6910076974402

$ws = $wb.Worksheets.Item("5")
$ws.Range("D2").Value = 0.08166881313964844
$ws.Range("E2").Value = 138.89541
$ws.Range("B3").Value = -963.1399746398824
$ws.Range("C3").Value = 0.09222430419959636
$ws.Range("D3").Value = 1016.0596363792858

$ws = $wb.Worksheets.Item("6")
$ws.Range("D2").Value = 0.11606455872436523
$ws.Range("E2").Value = 136.79724
$ws.Range("B3").Value = -953.069846438852
$ws.Range("C3").Value = 0.0
$ws.Range("D3").Value = 510.0949633082174

$ws = $wb.Worksheets.Item("7")
$ws.Range("D2").Value = 0.10300813076916504
$ws.Range("E2").Value = 136.77823
$ws.Range("B3").Value = -965.4153674106826
$ws.Range("C3").Value = 0.006156580694377443
$ws.Range("D3").Value = 947.7586012058363

$ws = $wb.Worksheets.Item("8")
$ws.Range("D2").Value = 0.11353465360314942
$ws.Range("E2").Value = 139.2123
$ws.Range("B3").Value = -954.6051029237938
$ws.Range("C3").Value = 0.029341303106506938
$ws.Range("D3").Value = 1209.368990072879

$ws = $wb.Worksheets.Item("9")
$ws.Range("D2").Value = 0.06519662681298828
$ws.Range("E2").Value = 138.57865
$ws.Range("B3").Value = -960.9871769768399
$ws.Range("C3").Value = 0.004023830828155816
$ws.Range("D3").Value = 416.2713530777032

$ws = $wb.Worksheets.Item("10")
$ws.Range("D2").Value = 0.11287795769885253
$ws.Range("E2").Value = 134.74132
$ws.Range("B3").Value = -965.36122253593
$ws.Range("C3").Value = 0.08386240805235516
$ws.Range("D3").Value = 946.3168402140044

Write-Host "Edit applied successfully"
